$wb = $excel.ActiveWorkbook

$wsIn = $wb.Worksheets.Item("ProductLoanInput")
$wsOut = $wb.Worksheets.Item("ProductLoanOutput")

$wsIn.Range("B1").Value = "4273-MS-EI-DB-DL-REC-RNI-INT-FFConMONTHLYonLASTSUNDAY-FIFC-1-FFROP-DAILY-FIFR-1-MD--1st"
$wsIn.Range("B2").Value = "427n"

$wsOut.Range("B1").Value = "4273-MS-EI-DB-DL-REC-RNI-INT-FFConMONTHLYonLASTSUNDAY-FIFC-1-FFROP-DAILY-FIFR-1-MD--1st"

$wsIn.Range("A27").Select()
$wsOut.Activate()
